# Updates the cryptocurrency price/volume snapshot in the active sheet.
# Mirrors a refreshed data pull: per-row Price (D) / Volume(1h) (E) values
# are updated, and the BabyDogeCoin/Quant rows (50-51) swap rank order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.866.93"
$ws.Range("E2").Value = "  +1.16%  "

$ws.Range("D3").Value = "1.622.65"
$ws.Range("E3").Value = "  +1.15%  "

$ws.Range("D4").Value = "'0.994"
$ws.Range("E4").Value = "  -0.45%  "

$ws.Range("D5").Value = "'213.57"
$ws.Range("E5").Value = "  +0.48%  "

$ws.Range("D6").Value = "'0.519"
$ws.Range("E6").Value = "  -0.77%  "

$ws.Range("E7").Value = "  -0.40%  "

$ws.Range("D8").Value = "'29.53"
$ws.Range("E8").Value = "  +9.75%  "

$ws.Range("E9").Value = "  +2.90%  "

$ws.Range("D10").Value = "'0.0607"
$ws.Range("E10").Value = "  +0.92%  "

$ws.Range("D11").Value = "'0.0911"
$ws.Range("E11").Value = "  +0.12%  "

$ws.Range("D12").Value = "1.854.94"
$ws.Range("E12").Value = "  +1.15%  "

$ws.Range("D13").Value = "1.616.34"
$ws.Range("E13").Value = "  +0.82%  "

$ws.Range("D14").Value = "'0.568"
$ws.Range("E14").Value = "  +5.70%  "

$ws.Range("E15").Value = "  +5.33%  "

$ws.Range("D16").Value = "29.869.30"
$ws.Range("E16").Value = "  +1.08%  "

$ws.Range("D17").Value = "'8.81"
$ws.Range("E17").Value = "  +15.92%  "

$ws.Range("D18").Value = "'64.36"
$ws.Range("E18").Value = "  +1.37%  "

$ws.Range("D19").Value = "'242.74"
$ws.Range("E19").Value = "  +0.93%  "

$ws.Range("D20").Value = "0.0₃0707"
$ws.Range("E20").Value = "  +2.27%  "

$ws.Range("E21").Value = "  -0.29%  "

$ws.Range("D22").Value = "'4.11"
$ws.Range("E22").Value = "  +3.32%  "

$ws.Range("D23").Value = "'9.58"
$ws.Range("E23").Value = "  +3.96%  "

$ws.Range("D24").Value = "'2.13"
$ws.Range("E24").Value = "  +2.21%  "

$ws.Range("D25").Value = "'156.64"
$ws.Range("E25").Value = "  +1.42%  "

$ws.Range("D26").Value = "'15.61"
$ws.Range("E26").Value = "  +2.33%  "

$ws.Range("E27").Value = "  +1.32%  "

$ws.Range("E28").Value = "  +2.99%  "

$ws.Range("D29").Value = "'0.995"
$ws.Range("E29").Value = "  -0.41%  "

$ws.Range("D30").Value = "'0.0488"
$ws.Range("E30").Value = "  +3.29%  "

$ws.Range("D31").Value = "'1.12"
$ws.Range("E31").Value = "  +5.41%  "

$ws.Range("E32").Value = "  +3.16%  "

$ws.Range("E33").Value = "  +3.89%  "

$ws.Range("D34").Value = "1.425.00"
$ws.Range("E34").Value = "  +1.16%  "

$ws.Range("E35").Value = "  +6.92%  "

$ws.Range("E36").Value = "  -0.64%  "

$ws.Range("E37").Value = "  +1.34%  "

$ws.Range("D38").Value = "'2.30"
$ws.Range("E38").Value = "  -0.36%  "

$ws.Range("D39").Value = "'0.0170"
$ws.Range("E39").Value = "  +2.76%  "

$ws.Range("D40").Value = "'0.555"
$ws.Range("E40").Value = "  +3.05%  "

$ws.Range("E41").Value = "  +3.27%  "

$ws.Range("D42").Value = "'1.98"
$ws.Range("E42").Value = "  +0.12%  "

$ws.Range("D43").Value = "'0.829"
$ws.Range("E43").Value = "  +3.88%  "

$ws.Range("D44").Value = "'54.08"
$ws.Range("E44").Value = "  +1.73%  "

$ws.Range("D45").Value = "'69.02"
$ws.Range("E45").Value = "  +4.94%  "

$ws.Range("E46").Value = "  +19.03%  "

$ws.Range("D47").Value = "'0.995"
$ws.Range("E47").Value = "  -0.32%  "

$ws.Range("E48").Value = "  +2.78%  "

$ws.Range("D49").Value = "1.763.00"
$ws.Range("E49").Value = "  +1.02%  "

# Quant moves up to row 50, BabyDogeCoin moves down to row 51
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "'88.19"
$ws.Range("E50").Value = "  +1.79%  "

$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0105"
$ws.Range("E51").Value = "  +3.30%  "
